# lesson 471 - Friday - Mexico catch up
#
# 1. Split the "Residence, sorroundings, ..." run so the corrected word
#    "surroundings" sits in its own run, and drop the _GoBack bookmark that
#    used to trail that paragraph (it moves to the new word-list paragraph
#    added in step 2).
# 2. Fill the first blank paragraph after the 2nd picture with
#    "On, , make, , , on, , , of, had" (with "make" spell-check-flagged),
#    and re-add the _GoBack bookmark at its end.
# 3. Fill the blank paragraph after the 3rd picture with the
#    "Analyst, shelved, materialized, ..." word list (each word flagged).

$d = $word.ActiveDocument

$para4Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00ED1AE5" w:rsidRPr="00ED1AE5" w:rsidRDefault="00ED1AE5"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00ED1AE5"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Residence, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>surroundings</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, solidity, architectural, manufacturers, specious, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">shelves, </w:t></w:r><w:r w:rsidR="00A661DE"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>infrastructure, inhabitable, residential</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para6Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">On, , </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>make</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, , , on, , , of, had</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para10Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Analyst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shelved</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>materialized</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>repaired</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spatial</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>residents</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inhabited</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>objections</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>densely</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>architectural</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 4: "Residence, sorroundings, ..." word list (before picture 2) ---
$d.Paragraphs(4).Range.InsertXML($para4Xml)

# --- Paragraph 6: first blank paragraph after picture 2 ---
$d.Paragraphs(6).Range.InsertXML($para6Xml)

# --- Paragraph 10: blank paragraph after picture 3 (last paragraph in the body) ---
# It is the very last paragraph in the document, so Word will not let its
# mark be replaced in place (a document always needs a trailing paragraph
# mark). Add a throwaway paragraph after it first, fill paragraph 10, then
# delete the throwaway paragraph.
$tailRng = $d.Range($d.Content.End, $d.Content.End)
$tailRng.InsertParagraphAfter()
$d.Paragraphs(10).Range.InsertXML($para10Xml)
$d.Paragraphs($d.Paragraphs.Count).Range.Delete()
